$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (row 2 totals, row 3 stats)
$ws.Range("F2").Value = 1.835123164180666
$ws.Range("I2").Value = 1.835123164180666
$ws.Range("B3").Value = 32
$ws.Range("C3").Value = 0.9234578654525802

# Add new rows of test-case data (videos et images)
$ws.Range("A4").Value = "https://www.amazon.fr/deals?ref_=nav_cs_gb"
$ws.Range("B4").Value = 66
$ws.Range("C4").Value = 1.278420568435639

$ws.Range("A5").Value = "https://www.amazon.fr/deal/3a51f27b?showVariations=true&pf_rd_r=JTA1B547JDKMBK143S96&pf_rd_t=Events&pf_rd_i=deals&pf_rd_p=3c3f3ff2-f80e-428b-aff0-b0531c852487&pf_rd_s=slot-14&ref=dlx_deals_gd_dcl_img_1_3a51f27b_dt_sl14_87"
$ws.Range("B5").Value = 69
$ws.Range("C5").Value = 1.599036199888587

$ws.Range("A6").Value = "https://www.amazon.fr/gp/your-account/order-history?ref_=ya_d_c_yo"
$ws.Range("B6").Value = 79
$ws.Range("C6").Value = 1.835123164180666
